# plot proportion disease, hosp, death
#
# 1) Rename Sheet1 -> SUMMARY
# 2) Add a new worksheet "other-papers" right after SUMMARY and populate it
# 3) Adjust SUMMARY's frozen-pane scroll position and activate other-papers

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original (only) sheet to SUMMARY
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Name = "SUMMARY"

# Scroll the frozen pane back up to the top (was topLeftCell A25)
$summary.Activate()
$summary.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$summary.Range("H42").Select()

# ---------------------------------------------------------------------
# 2. Add the new "other-papers" sheet right after SUMMARY
# ---------------------------------------------------------------------
$papers = $wb.Worksheets.Add($null, $summary)
$papers.Name = "other-papers"

# Column widths (A, B best-fit-like widths; E explicit width)
$papers.Columns.Item(1).ColumnWidth = 13.666666666666666
$papers.Columns.Item(2).ColumnWidth = 23.498697916666668
$papers.Columns.Item(5).ColumnWidth = 31.166666666666668

$xlCenter = -4108

function Set-Header($addr, $text) {
    $c = $papers.Range($addr)
    $c.Value2 = $text
    $c.HorizontalAlignment = $xlCenter
    $c.VerticalAlignment = $xlCenter
}

function Set-HeaderBlank($addr) {
    $c = $papers.Range($addr)
    $c.HorizontalAlignment = $xlCenter
    $c.VerticalAlignment = $xlCenter
}

function Set-Label($addr, $text) {
    $c = $papers.Range($addr)
    $c.Value2 = $text
    $c.VerticalAlignment = $xlCenter
}

function Set-Prop($addr, $value) {
    $c = $papers.Range($addr)
    $c.Value2 = $value
    $c.NumberFormat = "0.0000"
    $c.HorizontalAlignment = $xlCenter
    $c.VerticalAlignment = $xlCenter
}

function Set-Plain($addr, $value) {
    $c = $papers.Range($addr)
    $c.Value2 = $value
    $c.HorizontalAlignment = $xlCenter
    $c.VerticalAlignment = $xlCenter
}

function Set-Comment($addr, $text, [switch]$Center) {
    $c = $papers.Range($addr)
    $c.Value2 = $text
    $c.WrapText = $true
    if ($Center) {
        $c.HorizontalAlignment = $xlCenter
    }
}

# ---------------------------------------------------------------------
# Row 3: column headers (Low / High / comment)
# ---------------------------------------------------------------------
Set-HeaderBlank "A3"
Set-HeaderBlank "B3"
Set-Header "C3" "Low"
Set-Header "D3" "High"
Set-Plain  "E3" "comment"
$papers.Range("E3").HorizontalAlignment = $xlCenter
$papers.Range("E3").VerticalAlignment = $null

# ---------------------------------------------------------------------
# Biggerstaff 2015 block (rows 4-10)
# ---------------------------------------------------------------------
Set-Label "A4" "Biggerstaff 2015"
Set-Label "B4" "Final size prop"
Set-Prop  "C4" 0.2
Set-Prop  "D4" 0.3

Set-Label "A5" "Biggerstaff 2015"
Set-Label "B5" "Hospitalization prop"
Set-Prop  "C5" 0.005
Set-Prop  "D5" 0.042

Set-Label "A6" "Biggerstaff 2015"
Set-Label "B6" "Case fatality ratio"
Set-Prop  "C6" 0.0008
Set-Prop  "D6" 0.0053

Set-Label "A7" "Biggerstaff 2015"
Set-Label "B7" "Vax coverage"
Set-Prop  "C7" 0.8
Set-Prop  "D7" 0.8

Set-Label "A8" "Biggerstaff 2015"
Set-Label "B8" "Vax admin rate (10^6/week)"
Set-Prop  "C8" 10
Set-Prop  "D8" 30

Set-Label "A9" "Biggerstaff 2015"
Set-Label "B9" "Vax effectiveness"
Set-Prop  "C9" 0.43
Set-Prop  "D9" 0.8
$papers.Rows.Item(9).RowHeight = 30
Set-Comment "E9" "for > 60 yrs-old: 0.43 - 0.60`nfor <60yrs: 0.62-0.80"

Set-Label "A10" "Biggerstaff 2015"
Set-Label "B10" "Vax start (days)"
Set-Plain "C10" -112
Set-Plain "D10" 112

# ---------------------------------------------------------------------
# Greer 2015 block (rows 11-17)
# ---------------------------------------------------------------------
Set-Label "A11" "Greer 2015"
Set-Label "B11" "Final size prop"
Set-Prop  "C11" 0.2
Set-Prop  "D11" 0.5

Set-Label "A12" "Greer 2015"
Set-Label "B12" "Hospitalization prop"
Set-Prop  "C12" 0.01
Set-Prop  "D12" 0.01

Set-Label "A13" "Greer 2015"
Set-Label "B13" "Case fatality ratio"
Set-Prop  "C13" 0.004
Set-Prop  "D13" 0.004

Set-Label "A14" "Greer 2015"
Set-Label "B14" "Vax coverage"
Set-Prop  "C14" 0.26
Set-Prop  "D14" 0.75
$papers.Range("E14").Value2 = "age dependent"

Set-Label "A15" "Greer 2015"
Set-Label "B15" "Vax admin rate (10^6/week)"
Set-Prop  "C15" 0.15
Set-Prop  "D15" 3.75

Set-Label "A16" "Greer 2015"
Set-Label "B16" "Vax effectiveness"
Set-Prop  "C16" 0.3
Set-Prop  "D16" 0.7
$papers.Rows.Item(16).RowHeight = 30
Set-Comment "E16" "for >65 yrs-old: 0.3`nfor <65yrs: 0.7"

Set-Label "A17" "Greer 2015"
Set-Label "B17" "Vax start (days)"
Set-Plain "C17" 30
Set-Plain "D17" 90

# ---------------------------------------------------------------------
# Page setup to match the authored sheet
# ---------------------------------------------------------------------
$papers.PageSetup.LeftMargin = 0.75 * 72
$papers.PageSetup.RightMargin = 0.75 * 72
$papers.PageSetup.TopMargin = 1 * 72
$papers.PageSetup.BottomMargin = 1 * 72
$papers.PageSetup.HeaderMargin = 0.5 * 72
$papers.PageSetup.FooterMargin = 0.5 * 72
$papers.PageSetup.PaperSize = 9
$papers.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Zoom + selection + make other-papers the active/visible tab
# ---------------------------------------------------------------------
$papers.Activate()
$excel.ActiveWindow.Zoom = 205
$papers.Range("C18").Select()
